{"js": "// Replace each three-digit-by-one-digit multiplication expression in the\n// document with its new value (same position / order, one-to-one mapping).\n// All \"old\" values are unique in the document, so a plain matchCase search\n// for each one unambiguously finds its single run.\nconst replacements = [\n  { oldText: \"797\u00d77=\", newText: \"403\u00d74=\" },\n  { oldText: \"822\u00d76=\", newText: \"812\u00d78=\" },\n  { oldText: \"654\u00d75=\", newText: \"146\u00d79=\" },\n  { oldText: \"285\u00d79=\", newText: \"620\u00d75=\" },\n  { oldText: \"544\u00d72=\", newText: \"403\u00d77=\" },\n  { oldText: \"312\u00d72=\", newText: \"623\u00d74=\" },\n  { oldText: \"349\u00d78=\", newText: \"217\u00d76=\" },\n  { oldText: \"710\u00d72=\", newText: \"303\u00d78=\" },\n  { oldText: \"379\u00d73=\", newText: \"912\u00d78=\" },\n  { oldText: \"683\u00d77=\", newText: \"836\u00d73=\" },\n  { oldText: \"148\u00d77=\", newText: \"105\u00d74=\" },\n  { oldText: \"646\u00d76=\", newText: \"767\u00d72=\" },\n  { oldText: \"976\u00d79=\", newText: \"440\u00d72=\" },\n  { oldText: \"787\u00d74=\", newText: \"754\u00d79=\" },\n  { oldText: \"708\u00d74=\", newText: \"185\u00d72=\" },\n  { oldText: \"142\u00d79=\", newText: \"841\u00d75=\" },\n  { oldText: \"246\u00d76=\", newText: \"206\u00d79=\" },\n  { oldText: \"247\u00d77=\", newText: \"574\u00d75=\" },\n  { oldText: \"283\u00d73=\", newText: \"780\u00d74=\" },\n  { oldText: \"494\u00d79=\", newText: \"580\u00d73=\" },\n  { oldText: \"678\u00d74=\", newText: \"306\u00d78=\" },\n  { oldText: \"241\u00d73=\", newText: \"965\u00d73=\" },\n  { oldText: \"591\u00d74=\", newText: \"976\u00d73=\" },\n  { oldText: \"681\u00d78=\", newText: \"854\u00d73=\" },\n  { oldText: \"696\u00d76=\", newText: \"992\u00d72=\" },\n];\n\nconst body = context.document.body;\n\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication expression in the\n# document with its new value (same position / order, one-to-one mapping).\n# All \"old\" values are unique in the document, so a simple Find/Replace\n# (ReplaceAll) for each one unambiguously hits its single occurrence.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"797\u00d77=\"; New = \"403\u00d74=\" }\n    @{ Old = \"822\u00d76=\"; New = \"812\u00d78=\" }\n    @{ Old = \"654\u00d75=\"; New = \"146\u00d79=\" }\n    @{ Old = \"285\u00d79=\"; New = \"620\u00d75=\" }\n    @{ Old = \"544\u00d72=\"; New = \"403\u00d77=\" }\n    @{ Old = \"312\u00d72=\"; New = \"623\u00d74=\" }\n    @{ Old = \"349\u00d78=\"; New = \"217\u00d76=\" }\n    @{ Old = \"710\u00d72=\"; New = \"303\u00d78=\" }\n    @{ Old = \"379\u00d73=\"; New = \"912\u00d78=\" }\n    @{ Old = \"683\u00d77=\"; New = \"836\u00d73=\" }\n    @{ Old = \"148\u00d77=\"; New = \"105\u00d74=\" }\n    @{ Old = \"646\u00d76=\"; New = \"767\u00d72=\" }\n    @{ Old = \"976\u00d79=\"; New = \"440\u00d72=\" }\n    @{ Old = \"787\u00d74=\"; New = \"754\u00d79=\" }\n    @{ Old = \"708\u00d74=\"; New = \"185\u00d72=\" }\n    @{ Old = \"142\u00d79=\"; New = \"841\u00d75=\" }\n    @{ Old = \"246\u00d76=\"; New = \"206\u00d79=\" }\n    @{ Old = \"247\u00d77=\"; New = \"574\u00d75=\" }\n    @{ Old = \"283\u00d73=\"; New = \"780\u00d74=\" }\n    @{ Old = \"494\u00d79=\"; New = \"580\u00d73=\" }\n    @{ Old = \"678\u00d74=\"; New = \"306\u00d78=\" }\n    @{ Old = \"241\u00d73=\"; New = \"965\u00d73=\" }\n    @{ Old = \"591\u00d74=\"; New = \"976\u00d73=\" }\n    @{ Old = \"681\u00d78=\"; New = \"854\u00d73=\" }\n    @{ Old = \"696\u00d76=\"; New = \"992\u00d72=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($Null, $False, $False, $False, $False, $False, $True, 1, $False, $r.New, 2)\n}\n"}
